# Auto-generated edit script: updates crypto price/volume data and two row swaps
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '56.362.84'
$ws.Cells.Item(2, 5).Value = '  -0.56%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.324.06'
$ws.Cells.Item(3, 5).Value = '  -0.18%  '
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '511.17'
$ws.Cells.Item(5, 5).Value = '  -1.71%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '131.43'
$ws.Cells.Item(6, 5).Value = '  -2.36%  '
$ws.Cells.Item(7, 5).Value = '  +0.41%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.532'
$ws.Cells.Item(8, 5).Value = '  -1.18%  '
$ws.Cells.Item(9, 5).Value = '  -2.98%  '
$ws.Cells.Item(10, 5).Value = '  -0.32%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '5.26'
$ws.Cells.Item(11, 5).Value = '  +0.78%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.337'
$ws.Cells.Item(12, 5).Value = '  -1.47%  '
$ws.Cells.Item(13, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '2.741.73'
$ws.Cells.Item(13, 5).Value = '  -0.20%  '
$ws.Cells.Item(14, 2).Value = 'Avalanche'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '23.45'
$ws.Cells.Item(14, 5).Value = '  -1.02%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '56.351.58'
$ws.Cells.Item(15, 5).Value = '  -0.76%  '
$ws.Cells.Item(16, 5).Value = '  -1.89%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '2.330.80'
$ws.Cells.Item(17, 5).Value = '  -0.41%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '10.41'
$ws.Cells.Item(18, 5).Value = '  -0.59%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '325.57'
$ws.Cells.Item(19, 5).Value = '  +0.44%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '4.12'
$ws.Cells.Item(20, 5).Value = '  -2.58%  '
$ws.Cells.Item(21, 5).Value = '  +2.99%  '
$ws.Cells.Item(22, 5).Value = '  +0.02%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '61.44'
$ws.Cells.Item(23, 5).Value = '  +0.86%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '8.72'
$ws.Cells.Item(24, 5).Value = '  +11.05%  '
$ws.Cells.Item(25, 5).Value = '  -0.61%  '
$ws.Cells.Item(26, 5).Value = '  +0.41%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.29'
$ws.Cells.Item(27, 5).Value = '  +1.73%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '167.45'
$ws.Cells.Item(28, 5).Value = '  -1.88%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.66'
$ws.Cells.Item(29, 5).Value = '  -2.99%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.0₃0715'
$ws.Cells.Item(30, 5).Value = '  -4.09%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '6.09'
$ws.Cells.Item(31, 5).Value = '  -1.30%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '18.30'
$ws.Cells.Item(32, 5).Value = '  +0.01%  '
$ws.Cells.Item(33, 5).Value = '  +0.01%  '
$ws.Cells.Item(34, 5).Value = '  +0.68%  '
$ws.Cells.Item(35, 5).Value = '  +0.55%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.885'
$ws.Cells.Item(36, 5).Value = '  -3.86%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '3.90'
$ws.Cells.Item(37, 5).Value = '  -2.92%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '38.41'
$ws.Cells.Item(38, 5).Value = '  +1.39%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.54'
$ws.Cells.Item(39, 5).Value = '  -0.43%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '148.93'
$ws.Cells.Item(40, 5).Value = '  +8.63%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.54'
$ws.Cells.Item(42, 5).Value = '  -1.38%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '276.52'
$ws.Cells.Item(43, 5).Value = '  -0.67%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '5.02'
$ws.Cells.Item(44, 5).Value = '  -3.30%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0926'
$ws.Cells.Item(45, 5).Value = '  -0.88%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0492'
$ws.Cells.Item(46, 5).Value = '  -2.45%  '
$ws.Cells.Item(47, 5).Value = '  -1.23%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '18.00'
$ws.Cells.Item(48, 5).Value = '  +2.85%  '
$ws.Cells.Item(49, 2).Value = 'VeChain'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0213'
$ws.Cells.Item(49, 5).Value = '  -2.53%  '
$ws.Cells.Item(50, 2).Value = 'Polygon'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.375'
$ws.Cells.Item(50, 5).Value = '  -1.18%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '17.03'
$ws.Cells.Item(51, 5).Value = '  +0.83%  '
